$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.795.20"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.102.47"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.42"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.81"
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0842"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.75"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.413.61"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.00"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.801"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.094.27"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.810.72"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.60"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.03"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.88"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -5.82%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.64"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.14"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.35"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.57"
$ws.Range("E31").Value = "  +10.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.79"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("B35").Value = "THORChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.16"
$ws.Range("E35").Value = "  +11.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0615"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.99"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  +3.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.80"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.524.61"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("E44").Value = "  +7.39%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.77"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0911"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  +4.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.15"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.300.84"
$ws.Range("E51").Value = "  +0.24%  "
